# Apply the two changes described by the commit:
#  1. Slide 5's table switches from the custom "Table_0" style to the
#     built-in table style {837BDD73-7C98-4F33-B7F6-75A79B67788A}.
#  2. The presentation's theme swaps from the "Integral" (Red Violet)
#     colour scheme to the default "Office" colour scheme (the font
#     scheme and format scheme are identical between the two themes,
#     so only the 12 theme colours need to change).

$p = $ppt.ActivePresentation

# --- 1. Table style ---------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{837BDD73-7C98-4F33-B7F6-75A79B67788A}")

# --- 2. Theme colours ---------------------------------------------------
function ColorRef($r, $g, $b) {
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$colors = $p.SlideMaster.Theme.ThemeColorScheme

# Office theme colour scheme (replaces the Red Violet / Integral scheme)
$colors.Item(1).RGB  = ColorRef 0x00 0x00 0x00   # dk1
$colors.Item(2).RGB  = ColorRef 0xFF 0xFF 0xFF   # lt1
$colors.Item(3).RGB  = ColorRef 0x44 0x54 0x6A   # dk2
$colors.Item(4).RGB  = ColorRef 0xE7 0xE6 0xE6   # lt2
$colors.Item(5).RGB  = ColorRef 0x5B 0x9B 0xD5   # accent1
$colors.Item(6).RGB  = ColorRef 0xED 0x7D 0x31   # accent2
$colors.Item(7).RGB  = ColorRef 0xA5 0xA5 0xA5   # accent3
$colors.Item(8).RGB  = ColorRef 0xFF 0xC0 0x00   # accent4
$colors.Item(9).RGB  = ColorRef 0x44 0x72 0xC4   # accent5
$colors.Item(10).RGB = ColorRef 0x70 0xAD 0x47   # accent6
$colors.Item(11).RGB = ColorRef 0x05 0x63 0xC1   # hlink
$colors.Item(12).RGB = ColorRef 0x95 0x4F 0x72   # folHlink
